$d = $word.ActiveDocument

# 1) "Thread safety (only one instance of ApiClient)" ->
#    "Thread safety (only one instance of ApiClient which wraps HttpClient whose methods are thread safe)"
$d.Content.Find.Execute(
    "Thread safety (only one instance of ApiClient)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Thread safety (only one instance of ApiClient which wraps HttpClient whose methods are thread safe)", 2)

# 2) "Achitecture is set in that way that it can be scalable (separation of concerns)" ->
#    "Achitecture and programing stlye is set in that way that it can be scalable"
$d.Content.Find.Execute(
    "Achitecture is set in that way that it can be scalable (separation of concerns)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Achitecture and programing stlye is set in that way that it can be scalable", 2)

# 3) "documented (Summaries are everywhere)" -> "documented (Summaries are everywhere in code)"
$d.Content.Find.Execute(
    "documented (Summaries are everywhere)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "documented (Summaries are everywhere in code)", 2)

# 4) Move the "_GoBack" bookmark from after "- Implemented 2 days before deadline"
#    to right after the "- Mostly unit tested" paragraph text (re-adding with the
#    same name automatically replaces the previous location).
#    Note: a range collapsed exactly at (paragraph.End - 1) is mishandled by this
#    runtime, so we use a tiny range that straddles the paragraph mark
#    (End-1 .. End+1); Word still places both bookmark tags immediately after the
#    paragraph's text content.
$goBackPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Mostly unit tested*") {
        $goBackPara = $para
    }
}
$r = $goBackPara.Range
$target = $d.Range($r.End - 1, $r.End + 1)
$d.Bookmarks.Add("_GoBack", $target)
